# Adds new FB/Twitter sentiment columns (I:N) and two new recommendation
# notes (C10, C12) to the "Recommendations for Low SES Artists" sheet.
#
# Cell writes are ordered to match the shared-string table build order of
# the authored workbook: the I8:N8 headers first, then the H:N sentiment
# table (top to bottom), and only afterwards the two new column-C notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8: new header cells for the sentiment-percentage table (I8:N8)
# ---------------------------------------------------------------------
$ws.Range("I8").Value = "% Postive FB Posts"
$ws.Range("J8").Value = "% Neutral FB Posts"
$ws.Range("K8").Value = "% Negative FB Posts"
$ws.Range("L8").Value = "% Positive Tweets"
$ws.Range("M8").Value = "% Neutral Tweets"
$ws.Range("N8").Value = "% Negative Tweets"
$ws.Range("I8:N8").Font.Bold = $true
$ws.Range("I8:N8").WrapText = $true

# ---------------------------------------------------------------------
# Sentiment data table (H9:N16) - one row per artist
# ---------------------------------------------------------------------
$ws.Range("H9").Value = "Diego Torres"
$ws.Range("I9").Value = 35.8
$ws.Range("J9").Value = 60.5
$ws.Range("K9").Value = 3.7
$ws.Range("L9").Value = 57.5
$ws.Range("M9").Value = 42.5
$ws.Range("N9").Value = 0

$ws.Range("H10").Value = "Twenty One Pilots"
$ws.Range("I10").Value = 83.3
$ws.Range("J10").Value = 16.7
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 64.3
$ws.Range("M10").Value = 35.7
$ws.Range("N10").Value = 0

$ws.Range("H11").Value = "The Neighbourhood"
$ws.Range("I11").Value = 55.6
$ws.Range("J11").Value = 44.4
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 50
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 0

$ws.Range("H12").Value = "Harry Styles"
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 87.5
$ws.Range("M12").Value = 11.7
$ws.Range("N12").Value = 0.008

$ws.Range("H14").Value = "Calvin Harris"
$ws.Range("I14").Value = 64.5
$ws.Range("J14").Value = 35.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 76.9
$ws.Range("M14").Value = 19
$ws.Range("N14").Value = 4.1

$ws.Range("H15").Value = "Kygo"
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 77.8
$ws.Range("M15").Value = 22.2
$ws.Range("N15").Value = 0

$ws.Range("H16").Value = "Martin Garrix"
$ws.Range("I16").Value = 80
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 66
$ws.Range("M16").Value = 30.9
$ws.Range("N16").Value = 3.1

# Wrap text for the whole new block, matching the rest of the sheet
$ws.Range("H9:N16").WrapText = $true

# ---------------------------------------------------------------------
# New recommendation text in column C (added after the data table)
# ---------------------------------------------------------------------
$ws.Range("C10").Value = "Post more positive and less neutral posts on Facebook as well as Twitter"
$ws.Range("C12").Value = "Post fewer negative tweets"

# ---------------------------------------------------------------------
# View state: Excel re-froze the pane at the top of the data and moved
# the active selection down one row after the new row 8 headers were
# inserted.
# ---------------------------------------------------------------------
$ws.Range("B2").Select() | Out-Null
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("C12").Select() | Out-Null
